# "able to insert data to cell and then increment by one. Each scan will be
#  add to a cell in Excel"
#
# The sheet used to hold an incrementing counter in column A (1..9) plus a
# repeated scan-id string in column D. The workbook now just appends each
# scanned id straight into column A, one row per scan, starting at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old layout (counter column A + shared-string column D).
$ws.Cells.Clear()

# Each barcode/QR "scan" gets written to the next free cell in column A.
$scanId = "['8a9ab3409000']"
$startRow = 2
$scanCount = 5

$endRow = $startRow + $scanCount - 1
for ($i = 0; $i -lt $scanCount; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $scanId

    # Touch the built-in Hyperlink / Followed Hyperlink cell styles (the
    # scan cells are rendered as links to their source), then drop back to
    # Normal so the data itself stays plain text.
    $ws.Hyperlinks.Add($cell, "", "", "", $scanId) | Out-Null
}
$dataRange = "A" + $startRow + ":A" + $endRow
$ws.Range("A" + $startRow).Style = "Followed Hyperlink"
$ws.Range($dataRange).Style = "Normal"
$ws.Hyperlinks.Delete()

$ws.PageSetup.Orientation = 1

# Leave the cursor parked below the data, ready for the next scan.
$ws.Range("B9").Select() | Out-Null
